# The edit rotates the full content (every column, A through AY) of 20 data
# rows (rows 2-20 and 22) around a single 20-cycle. Row 21, the header row 1,
# and row 23 are left untouched.
#
# Cycle (dest <- src, i.e. dest row ends up holding src row's original
# content):
#   2<-16, 16<-17, 17<-8, 8<-14, 14<-7, 7<-5, 5<-13, 13<-15, 15<-20, 20<-10,
#   10<-3, 3<-18, 18<-19, 19<-12, 12<-22, 22<-11, 11<-4, 4<-9, 9<-6, 6<-2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycle = @(2, 16, 17, 8, 14, 7, 5, 13, 15, 20, 10, 3, 18, 19, 12, 22, 11, 4, 9, 6)
$n = $cycle.Length

# Use Copy / PasteSpecial (a true cell copy) rather than reading/writing
# through .Value — a plain .Value assignment lets Excel "helpfully" reinterpret
# text that looks like a date (e.g. the "2023-08-31" Startdatum/Slutdatum
# strings) as a real date serial, which would corrupt those cells. Copying
# the cell outright preserves its original type/format exactly.

# Step 1: stash every row involved in the cycle into scratch rows far below
# the real data, since it's a cycle and we can't safely overwrite a row
# before every other row has read what it needs from it.
$scratchBase = 1000
for ($i = 0; $i -lt $n; $i++) {
    $r = $cycle[$i]
    $ws.Range("A$r`:AY$r").Copy()
    $scratchRow = $scratchBase + $i
    $ws.Range("A$scratchRow`:AY$scratchRow").PasteSpecial()
}
$excel.CutCopyMode = $false

# Step 2: redistribute from the scratch rows back onto the real rows.
# Destination cycle[i] receives the content that originally lived at
# cycle[i+1] (wrapping), which is now safely parked at scratch row
# (scratchBase + i + 1).
for ($i = 0; $i -lt $n; $i++) {
    $destRow = $cycle[$i]
    $srcScratchRow = $scratchBase + (($i + 1) % $n)
    $ws.Range("A$srcScratchRow`:AY$srcScratchRow").Copy()
    $ws.Range("A$destRow`:AY$destRow").PasteSpecial()
}
$excel.CutCopyMode = $false

# Step 3: clean up the scratch rows.
for ($i = 0; $i -lt $n; $i++) {
    $scratchRow = $scratchBase + $i
    $ws.Range("A$scratchRow`:AY$scratchRow").Clear()
}
